$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 24.35712233333334
$ws.Range("H2").Value = 73.07136700000001
$ws.Range("I2").Value = 0.3750500562097488
$ws.Range("J2").Value = 0.3750500562097488
$ws.Range("M2").Value = 3.556762333333333
$ws.Range("N2").Value = 10.670287
$ws.Range("O2").Value = 0.04280930450251701
$ws.Range("P2").Value = 0.04280930450251701
$ws.Range("Q2").Value = 86.63249526359212
$ws.Range("R2").Value = 779.6924573723292
$ws.Range("S2").Value = 0.01605563205996926
$ws.Range("T2").Value = 0.01605563205996926
$ws.Range("G3").Value = 24.35712233333334
$ws.Range("H3").Value = 73.07136700000001
$ws.Range("I3").Value = 0.3750500562097488
$ws.Range("J3").Value = 0.3750500562097488
$ws.Range("O3").Value = 0.5686906263805706
$ws.Range("P3").Value = 0.5686906263805704
$ws.Range("Q3").Value = 1150.849997889298
$ws.Range("R3").Value = 10357.64998100368
$ws.Range("S3").Value = 0.2132874513899903
$ws.Range("T3").Value = 0.2132874513899902
$ws.Range("G4").Value = 24.35712233333334
$ws.Range("H4").Value = 73.07136700000001
$ws.Range("I4").Value = 0.3750500562097488
$ws.Range("J4").Value = 0.3750500562097488
$ws.Range("M4").Value = 24.53173066666666
$ws.Range("N4").Value = 73.595192
$ws.Range("O4").Value = 0.2952646900921413
$ws.Range("P4").Value = 0.2952646900921412
$ws.Range("Q4").Value = 597.5223648963849
$ws.Range("R4").Value = 5377.701284067464
$ws.Range("S4").Value = 0.1107390386158117
$ws.Range("T4").Value = 0.1107390386158116
$ws.Range("G5").Value = 24.35712233333334
$ws.Range("H5").Value = 73.07136700000001
$ws.Range("I5").Value = 0.3750500562097488
$ws.Range("J5").Value = 0.3750500562097488
$ws.Range("M5").Value = 7.746355333333334
$ws.Range("N5").Value = 23.239066
$ws.Range("O5").Value = 0.09323537902477132
$ws.Range("P5").Value = 0.0932353790247713
$ws.Range("Q5").Value = 188.6789244914691
$ws.Range("R5").Value = 1698.110320423222
$ws.Range("S5").Value = 0.03496793414397772
$ws.Range("T5").Value = 0.03496793414397771
$ws.Range("I6").Value = 0.2805618708302703
$ws.Range("J6").Value = 0.2805618708302702
$ws.Range("M6").Value = 3.556762333333333
$ws.Range("N6").Value = 10.670287
$ws.Range("O6").Value = 0.04280930450251701
$ws.Range("P6").Value = 0.04280930450251701
$ws.Range("Q6").Value = 64.80674924163934
$ws.Range("R6").Value = 583.260743174754
$ws.Range("S6").Value = 0.01201065856016889
$ws.Range("T6").Value = 0.01201065856016888
$ws.Range("I7").Value = 0.2805618708302703
$ws.Range("J7").Value = 0.2805618708302702
$ws.Range("O7").Value = 0.5686906263805706
$ws.Range("P7").Value = 0.5686906263805704
$ws.Range("S7").Value = 0.1595529060609711
$ws.Range("T7").Value = 0.1595529060609711
$ws.Range("I8").Value = 0.2805618708302703
$ws.Range("J8").Value = 0.2805618708302702
$ws.Range("M8").Value = 24.53173066666666
$ws.Range("N8").Value = 73.595192
$ws.Range("O8").Value = 0.2952646900921413
$ws.Range("P8").Value = 0.2952646900921412
$ws.Range("Q8").Value = 446.9856484023626
$ws.Range("R8").Value = 4022.870835621264
$ws.Range("S8").Value = 0.08284001384237114
$ws.Range("T8").Value = 0.08284001384237109
$ws.Range("I9").Value = 0.2805618708302703
$ws.Range("J9").Value = 0.2805618708302702
$ws.Range("M9").Value = 7.746355333333334
$ws.Range("N9").Value = 23.239066
$ws.Range("O9").Value = 0.09323537902477132
$ws.Range("P9").Value = 0.0932353790247713
$ws.Range("Q9").Value = 141.1441250710413
$ws.Range("R9").Value = 1270.297125639372
$ws.Range("S9").Value = 0.02615829236675918
$ws.Range("T9").Value = 0.02615829236675917
$ws.Range("G10").Value = 22.31748066666667
$ws.Range("H10").Value = 66.952442
$ws.Range("I10").Value = 0.3436437303202491
$ws.Range("J10").Value = 0.343643730320249
$ws.Range("M10").Value = 3.556762333333333
$ws.Range("N10").Value = 10.670287
$ws.Range("O10").Value = 0.04280930450251701
$ws.Range("P10").Value = 0.04280930450251701
$ws.Range("Q10").Value = 79.3779746100949
$ws.Range("R10").Value = 714.401771490854
$ws.Range("S10").Value = 0.01471114909166038
$ws.Range("T10").Value = 0.01471114909166038
$ws.Range("G11").Value = 22.31748066666667
$ws.Range("H11").Value = 66.952442
$ws.Range("I11").Value = 0.3436437303202491
$ws.Range("J11").Value = 0.343643730320249
$ws.Range("O11").Value = 0.5686906263805706
$ws.Range("P11").Value = 0.5686906263805704
$ws.Range("Q11").Value = 1054.478941585742
$ws.Range("R11").Value = 9490.310474271681
$ws.Range("S11").Value = 0.1954269682475783
$ws.Range("T11").Value = 0.1954269682475782
$ws.Range("G12").Value = 22.31748066666667
$ws.Range("H12").Value = 66.952442
$ws.Range("I12").Value = 0.3436437303202491
$ws.Range("J12").Value = 0.343643730320249
$ws.Range("M12").Value = 24.53173066666666
$ws.Range("N12").Value = 73.595192
$ws.Range("O12").Value = 0.2952646900921413
$ws.Range("P12").Value = 0.2952646900921412
$ws.Range("Q12").Value = 547.4864248732071
$ws.Range("R12").Value = 4927.377823858864
$ws.Range("S12").Value = 0.1014658595351157
$ws.Range("T12").Value = 0.1014658595351157
$ws.Range("G13").Value = 22.31748066666667
$ws.Range("H13").Value = 66.952442
$ws.Range("I13").Value = 0.3436437303202491
$ws.Range("J13").Value = 0.343643730320249
$ws.Range("M13").Value = 7.746355333333334
$ws.Range("N13").Value = 23.239066
$ws.Range("O13").Value = 0.09323537902477132
$ws.Range("P13").Value = 0.0932353790247713
$ws.Range("Q13").Value = 172.8791353887969
$ws.Range("R13").Value = 1555.912218499172
$ws.Range("S13").Value = 0.03203975344589472
$ws.Range("T13").Value = 0.03203975344589471
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.04834033333333334
$ws.Range("H14").Value = 0.145021
$ws.Range("I14").Value = 0.0007443426397318391
$ws.Range("J14").Value = 0.0007443426397318388
$ws.Range("M14").Value = 3.556762333333333
$ws.Range("N14").Value = 10.670287
$ws.Range("O14").Value = 0.04280930450251701
$ws.Range("P14").Value = 0.04280930450251701
$ws.Range("Q14").Value = 0.1719350767807778
$ws.Range("R14").Value = 1.547415691027
$ws.Range("S14").Value = 0.00003186479071848762
$ws.Range("T14").Value = 0.0000318647907184876
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.04834033333333334
$ws.Range("H15").Value = 0.145021
$ws.Range("I15").Value = 0.0007443426397318391
$ws.Range("J15").Value = 0.0007443426397318388
$ws.Range("O15").Value = 0.5686906263805706
$ws.Range("P15").Value = 0.5686906263805704
$ws.Range("Q15").Value = 2.284033054204444
$ws.Range("R15").Value = 20.55629748784
$ws.Range("S15").Value = 0.0004233006820308669
$ws.Range("T15").Value = 0.0004233006820308667
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.04834033333333334
$ws.Range("H16").Value = 0.145021
$ws.Range("I16").Value = 0.0007443426397318391
$ws.Range("J16").Value = 0.0007443426397318388
$ws.Range("M16").Value = 24.53173066666666
$ws.Range("N16").Value = 73.595192
$ws.Range("O16").Value = 0.2952646900921413
$ws.Range("P16").Value = 0.2952646900921412
$ws.Range("Q16").Value = 1.185872037670222
$ws.Range("R16").Value = 10.672848339032
$ws.Range("S16").Value = 0.0002197780988427878
$ws.Range("T16").Value = 0.0002197780988427877
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.04834033333333334
$ws.Range("H17").Value = 0.145021
$ws.Range("I17").Value = 0.0007443426397318391
$ws.Range("J17").Value = 0.0007443426397318388
$ws.Range("M17").Value = 7.746355333333334
$ws.Range("N17").Value = 23.239066
$ws.Range("O17").Value = 0.09323537902477132
$ws.Range("P17").Value = 0.0932353790247713
$ws.Range("Q17").Value = 0.3744613989317778
$ws.Range("R17").Value = 3.370152590386001
$ws.Range("S17").Value = 0.00006939906813969682
$ws.Range("T17").Value = 0.00006939906813969678
